$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp on row 10 (slightly adjusted value for same date)
$ws.Range("A10").Value = 44323.77488254398

# Add the new row 11 with the newly retrieved data
$ws.Range("A11").Value = 44324.77066125402
$ws.Range("B11").Value = 74182
$ws.Range("C11").Value = 62392
$ws.Range("D11").Value = 3192
$ws.Range("E11").Value = 2052
$ws.Range("F11").Value = 1446
$ws.Range("G11").Value = 19313
$ws.Range("H11").Value = 1356
$ws.Range("I11").Value = 835
$ws.Range("J11").Value = 216

# A11 should carry the same date/time number format as the rest of column A
$ws.Range("A11").NumberFormat = $ws.Range("A10").NumberFormat
